$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 367.66666
$ws.Range("I2").Value = 367.66666
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 367.66666
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -254.66666
$ws.Range("N2").ClearContents()

$ws.Range("H112").Value = 1622.7778
$ws.Range("I112").Value = 1434.25
$ws.Range("J112").Value = 1676.6428
$ws.Range("K112").Value = 4302.75
$ws.Range("L112").Value = 5029.928400000001
$ws.Range("M112").Value = -3194.75
$ws.Range("N112").Value = -7245.928400000001

$ws.Range("H138").Value = 4133.697
$ws.Range("I138").Value = 5417.6665
$ws.Range("J138").Value = 3848.3704
$ws.Range("K138").Value = 16252.9995
$ws.Range("L138").Value = 11545.1112
$ws.Range("M138").Value = -11112.9995
$ws.Range("N138").Value = -21825.1112

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 405.33334
$ws.Range("I22").Value = 16
$ws.Range("J22").Value = 600
$ws.Range("K22").Value = 16
$ws.Range("L22").Value = 600
$ws.Range("M22").Value = 283
$ws.Range("N22").Value = -1198

$ws.Range("H25").Value = 800
$ws.Range("I25").Value = 800
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 800
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -398
$ws.Range("N25").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1587.3334
$ws.Range("I20").Value = 1447.8334
$ws.Range("J20").Value = 1866.3334
$ws.Range("K20").Value = 1447.8334
$ws.Range("L20").Value = 1866.3334
$ws.Range("M20").Value = -1200.8334
$ws.Range("N20").Value = -2360.3334

$ws.Range("H24").Value = 4432.3335
$ws.Range("I24").Value = 4432.3335
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 4432.3335
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -4197.3335
$ws.Range("N24").ClearContents()

$ws.Range("H29").Value = 7099.6
$ws.Range("I29").Value = 7099.6
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 7099.6
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -6810.6

$ws.Range("H35").Value = 30000
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 30000
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 30000
$ws.Range("N35").Value = -30620

$ws.Range("H99").Value = 1831.4546
$ws.Range("I99").Value = 1764.6
$ws.Range("J99").Value = 2500
$ws.Range("K99").Value = 1764.6
$ws.Range("L99").Value = 2500
$ws.Range("M99").Value = -266.5999999999999
$ws.Range("N99").Value = -5496

$ws.Range("H134").Value = 1859.2128
$ws.Range("I134").Value = 1137.742
$ws.Range("J134").Value = 3257.0625
$ws.Range("K134").Value = 3413.226
$ws.Range("L134").Value = 9771.1875
$ws.Range("M134").Value = -878.2259999999997
$ws.Range("N134").Value = -14841.1875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 46686.31
$ws.Range("I31").Value = 3990.2727
$ws.Range("J31").Value = 77996.734
$ws.Range("K31").Value = 3990.2727
$ws.Range("L31").Value = 77996.734
$ws.Range("M31").Value = -3695.2727
$ws.Range("N31").Value = -78586.734

$ws.Range("H34").Value = 46686.31
$ws.Range("I34").Value = 3990.2727
$ws.Range("J34").Value = 77996.734
$ws.Range("K34").Value = 3990.2727
$ws.Range("L34").Value = 77996.734
$ws.Range("M34").Value = -3788.2727
$ws.Range("N34").Value = -78400.734

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 659.1875
$ws.Range("I2").Value = 69.125
$ws.Range("J2").Value = 1249.25
$ws.Range("K2").Value = 69.125
$ws.Range("L2").Value = 1249.25
$ws.Range("M2").Value = 43.875
$ws.Range("N2").Value = -1475.25

$ws.Range("H20").Value = 32003.334
$ws.Range("I20").Value = 32003.334
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 32003.334
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -31758.334
$ws.Range("N20").ClearContents()

$ws.Range("H24").Value = 29335.334
$ws.Range("I24").Value = 29335.334
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 29335.334
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -29162.334
$ws.Range("N24").ClearContents()

$ws.Range("H70").Value = 8378.799999999999
$ws.Range("I70").Value = 8160.2
$ws.Range("J70").Value = 8597.4
$ws.Range("K70").Value = 8160.2
$ws.Range("L70").Value = 8597.4
$ws.Range("M70").Value = -7890.2
$ws.Range("N70").Value = -9137.4

$ws.Range("H73").Value = 8378.799999999999
$ws.Range("I73").Value = 8160.2
$ws.Range("J73").Value = 8597.4
$ws.Range("K73").Value = 8160.2
$ws.Range("L73").Value = 8597.4
$ws.Range("M73").Value = -7224.2
$ws.Range("N73").Value = -10469.4

$ws.Range("H80").Value = 8500.125
$ws.Range("I80").Value = 7499.5
$ws.Range("J80").Value = 9500.75
$ws.Range("K80").Value = 7499.5
$ws.Range("L80").Value = 9500.75
$ws.Range("M80").Value = -6501.5
$ws.Range("N80").Value = -11496.75

$ws.Range("H83").Value = 8500.125
$ws.Range("I83").Value = 7499.5
$ws.Range("J83").Value = 9500.75
$ws.Range("K83").Value = 37497.5
$ws.Range("L83").Value = 47503.75
$ws.Range("M83").Value = -32505.5
$ws.Range("N83").Value = -57487.75

$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 170

$ws.Range("H132").Value = 22640.117
$ws.Range("I132").Value = 34308.934
$ws.Range("J132").Value = 4553.45
$ws.Range("K132").Value = 102926.802
$ws.Range("L132").Value = 13660.35
$ws.Range("M132").Value = -100396.802
$ws.Range("N132").Value = -18720.35

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 12482.333
$ws.Range("I22").Value = 4999
$ws.Range("J22").Value = 13979
$ws.Range("K22").Value = 4999
$ws.Range("L22").Value = 13979
$ws.Range("M22").Value = -4704
$ws.Range("N22").Value = -14569

$ws.Range("H27").Value = 12482.333
$ws.Range("I27").Value = 4999
$ws.Range("J27").Value = 13979
$ws.Range("K27").Value = 4999
$ws.Range("L27").Value = 13979
$ws.Range("M27").Value = -4892
$ws.Range("N27").Value = -14193

$ws.Range("H82").Value = 2583.375
$ws.Range("I82").Value = 1646.7
$ws.Range("J82").Value = 4144.5
$ws.Range("K82").Value = 1646.7
$ws.Range("L82").Value = 4144.5
$ws.Range("M82").Value = -1285.7
$ws.Range("N82").Value = -4866.5

$ws.Range("H85").Value = 2583.375
$ws.Range("I85").Value = 1646.7
$ws.Range("J85").Value = 4144.5
$ws.Range("K85").Value = 1646.7
$ws.Range("L85").Value = 4144.5
$ws.Range("M85").Value = -398.7
$ws.Range("N85").Value = -6640.5

$ws.Range("H100").Value = 2946.5386
$ws.Range("I100").Value = 1970.45
$ws.Range("J100").Value = 6200.1665
$ws.Range("K100").Value = 1970.45
$ws.Range("L100").Value = 6200.1665
$ws.Range("M100").Value = -1429.45
$ws.Range("N100").Value = -7282.1665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").ClearContents()

$ws.Range("H31").Value = 15333.333
$ws.Range("I31").Value = 15333.333
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 15333.333
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -14985.333
$ws.Range("N31").ClearContents()

$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()

$ws.Range("H51").Value = 24679.092
$ws.Range("I51").Value = 17470
$ws.Range("J51").Value = 25400
$ws.Range("K51").Value = 17470
$ws.Range("L51").Value = 25400
$ws.Range("M51").Value = -16960
$ws.Range("N51").Value = -26420

$ws.Range("H52").Value = 10909.091
$ws.Range("I52").Value = 15000
$ws.Range("J52").Value = 10500
$ws.Range("K52").Value = 15000
$ws.Range("L52").Value = 10500
$ws.Range("M52").Value = -14774
$ws.Range("N52").Value = -10952
